# Error Calculations and Plots
# Apply the data corrections to the missing_data worksheet:
#  - two records (RM 232, SC 92) are removed from the table entirely
#  - several cells have their "missing" (blank) / "present" (numeric) status
#    corrected for columns C (D col-letter header) and D (E col-letter header)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two rows that no longer belong in the cleaned data set.
#    "RM 232" is currently row 26; deleting it shifts everything below up by one.
$ws.Rows(26).Delete()
#    "SC 92" is now (after the shift) row 27; deleting it shifts everything below up again.
$ws.Rows(27).Delete()

# 2) Fix up individual cell values (missing <-> present) to match the corrected data.
$ws.Range("E5").Value = ""
$ws.Range("D6").Value = -14.2
$ws.Range("D8").Value = ""
$ws.Range("E11").Value = -7.9
$ws.Range("D19").Value = -15.5
$ws.Range("E19").Value = ""
$ws.Range("D21").Value = ""
$ws.Range("D23").Value = -13.9
$ws.Range("E23").Value = -7
$ws.Range("E25").Value = -7.1
$ws.Range("B26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = ""
$ws.Range("E30").Value = -5.7
$ws.Range("E33").Value = -10.7
